{"js": "// The placeholder \"{{forma_captura}}\" must become \"{{forma_capturado}}\".\n// This is done by appending a brand-new run containing \"do\" right after the\n// existing \"forma_captura\" run (inside the same spell-check span), rather\n// than rewriting the text of the existing run in place.\n\nconst body = context.document.body;\n\n// Find the \"forma_captura\" text.\nconst results = body.search(\"forma_captura\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n    throw new Error(\"Could not find 'forma_captura' in the document\");\n}\n\nconst match = results.items[0];\n\n// Collapse to the caret position right after \"forma_captura\" (and before\n// the closing \"}}\") and insert the new text there as its own range.\nconst caret = match.getRange(\"End\");\nconst inserted = caret.insertText(\"do\", \"Replace\");\nawait context.sync();\n\n// Nudge the formatting so the engine keeps \"do\" as a separate run instead of\n// silently merging it back into the preceding \"forma_captura\" run (both runs\n// end up with identical Nunito/32pt formatting either way).\ninserted.font.bold = true;\nawait context.sync();\ninserted.font.bold = false;\nawait context.sync();\n", "ps1": "# The placeholder \"{{forma_captura}}\" must become \"{{forma_capturado}}\".\n# This is done by appending a brand-new run containing \"do\" right after the\n# existing \"forma_captura\" run (inside the same spell-check span), rather\n# than rewriting the text of the existing run in place.\n\n$d = $word.ActiveDocument\n\n# Locate the \"forma_captura\" text and collapse the found range to its end\n# (the caret position immediately after the \"a\" in \"forma_captura\", and\n# before the closing \"}}\").\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"forma_captura\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find 'forma_captura' in the document\"\n}\n\n$rng.Collapse(0)  # wdCollapseEnd\n\n# Insert the new text as its own run right at the caret.\n$rng.InsertAfter(\"do\")\n\n# Nudge the formatting so the engine keeps \"do\" as a separate run instead of\n# silently merging it back into the preceding \"forma_captura\" run (both runs\n# end up with identical Nunito/32pt formatting either way).\n$rng.Bold = 1\n$rng.Bold = 0\n"}
